# RippleTemplate_Combination.xlsx edit:
#  - add a new "Assay" worksheet (settings/value table) at the end of the workbook
#  - make "Patterns" the active sheet/tab again, with O10 selected
#  - (Compounds loses its tabSelected flag automatically once another sheet becomes active)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Assay" sheet after the last existing sheet (Barcodes)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

# ---------------------------------------------------------------------------
# 2. Populate the settings/value table
# ---------------------------------------------------------------------------
$assay.Range("A1").Value = "Setting"
$assay.Range("B1").Value = "Value"

$assay.Range("A2").Value = "DMSO Tolerance"
$assay.Range("B2").Value = 0.005

$assay.Range("A3").Value = "Well Volume (µL)"
$assay.Range("B3").Value = 25

$assay.Range("A4").Value = "Backfill (µL)"
$assay.Range("B4").Value = 10

$assay.Range("A5").Value = "Allowed Error"
$assay.Range("B5").Value = 0.1

$assay.Range("A6").Value = "Destination Replicates"
$assay.Range("B6").Value = 1

$assay.Range("A7").Value = "Use Intermediate Plates"
$assay.Range("B7").Value = 1

$assay.Range("A8").Value = "DMSO Normalization"
$assay.Range("B8").Value = 1

# Leave the Assay sheet's own stored selection on its full data block
$assay.Range("A1:B8").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-activate "Patterns" as the visible/selected tab, cursor on O10
# ---------------------------------------------------------------------------
$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate()
$patterns.Range("O10").Select() | Out-Null
